$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows before row 24, shifting existing rows 24-40 down to 30-46
$ws.Rows("24:29").Insert()

# First block: rows 24-26 (OVI Alcohol / Drugs 1st-3rd, statute 4511.19(A)(1))
$ws.Range("A24").Value = "OVI Alcohol / Drugs 1st"
$ws.Range("A25").Value = "OVI Alcohol / Drugs 2nd"
$ws.Range("A26").Value = "OVI Alcohol / Drugs 3rd "
$ws.Range("B24").Value = "4511.19(A)(1)*"
$ws.Range("B25").Value = "4511.19(A)(1)**"
$ws.Range("B26").Value = "4511.19(A)(1)***"

# Second block: rows 27-29 (OVI Refusal 1st-3rd, statute 4511.19(A)(2))
$ws.Range("A27").Value = "OVI Refusal 1st / Prior in Past 20 Years"
$ws.Range("A28").Value = "OVI Refusal 2nd / Prior in Past 20 Years"
$ws.Range("A29").Value = "OVI Refusal 3rd / Prior in Past 20 Years"
$ws.Range("B27").Value = "4511.19(A)(2)*"
$ws.Range("B28").Value = "4511.19(A)(2)**"
$ws.Range("B29").Value = "4511.19(A)(2)***"

# Column C (degree)
$ws.Range("C24").Value = "M1"
$ws.Range("C25").Value = "M1"
$ws.Range("C26").Value = "UCM"
$ws.Range("C27").Value = "M1"
$ws.Range("C28").Value = "M1"
$ws.Range("C29").Value = "UCM"

# Column D (offense type)
$ws.Range("D24").Value = "Moving"
$ws.Range("D25").Value = "Moving"
$ws.Range("D26").Value = "Moving"
$ws.Range("D27").Value = "Moving"
$ws.Range("D28").Value = "Moving"
$ws.Range("D29").Value = "Moving"

# Match original formatting/style (text number format) of column B for the new rows
$ws.Range("B24:B29").NumberFormat = "@"

$ws.Columns("B").ColumnWidth = 18.7

$ws.Range("D30").Select()
